$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 3 ("Background"): tighten the wording of one of the bullet
# questions in the Content Placeholder.
# -----------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
$q = $body3.Paragraphs(3)
$q.Runs(1).Text = "How correlated are the features of a house to the sale price of the house?"

# -----------------------------------------------------------------
# Slide 7 ("Preparing the data"): rework the first bullet about
# categorical encoding and add follow-up sub-bullets / a new
# sub-bullet under the train/validation split bullet.
# -----------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$body7 = $s7.Shapes.Item(2).TextFrame.TextRange

# Replace the text of the first paragraph's run in place (keeps it a
# single run rather than splitting it).
$first = $body7.Paragraphs(1)
$first.Runs(1).Text = "Categorical" + [char]0x2013 + "nominal features will be one-hot encoded. "

# Insert two new sub-bullets right after it.
$first.InsertAfter("`rSome are quite large (city/zip code)`rWill we need to approach this another way?")
$body7.Paragraphs(2).IndentLevel = 2
$body7.Paragraphs(3).IndentLevel = 3

# Insert a new sub-bullet after the "dataset will be randomized..." bullet
# (now the 4th paragraph).
$split = $body7.Paragraphs(4)
$split.InsertAfter("`rWill we need to narrow the data down to one state, city, or zip code?")
$body7.Paragraphs(5).IndentLevel = 2
